$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1824
$ws.Cells.Item(3, 6).Value = 19
$ws.Cells.Item(6, 6).Value = 1072
$ws.Cells.Item(9, 6).Value = 561
$ws.Cells.Item(10, 6).Value = 47
$ws.Cells.Item(11, 6).Value = 443
$ws.Cells.Item(12, 6).Value = 196
$ws.Cells.Item(13, 6).Value = 1347
$ws.Cells.Item(14, 6).Value = 1187
$ws.Cells.Item(15, 6).Value = 1386
$ws.Cells.Item(16, 6).Value = 10
$ws.Cells.Item(17, 6).Value = 86
$ws.Cells.Item(18, 6).Value = 273
$ws.Cells.Item(19, 6).Value = 1533
$ws.Cells.Item(21, 6).Value = 762
$ws.Cells.Item(22, 6).Value = 303
$ws.Cells.Item(24, 6).Value = 104
$ws.Cells.Item(25, 6).Value = 1141
$ws.Cells.Item(26, 6).Value = 306
$ws.Cells.Item(27, 6).Value = 6
$ws.Cells.Item(30, 6).Value = 986
$ws.Cells.Item(31, 6).Value = 208383
$ws.Cells.Item(32, 6).Value = 927
$ws.Cells.Item(35, 6).Value = 873
$ws.Cells.Item(37, 6).Value = 11
$ws.Cells.Item(38, 6).Value = 814
$ws.Cells.Item(39, 6).Value = 1532
$ws.Cells.Item(40, 6).Value = 82
$ws.Cells.Item(41, 6).Value = 12
$ws.Cells.Item(42, 6).Value = 771
$ws.Cells.Item(44, 6).Value = 756
$ws.Cells.Item(45, 6).Value = 99

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 105
$ws.Cells.Item(6, 6).Value = 135
$ws.Cells.Item(8, 6).Value = 81
$ws.Cells.Item(11, 6).Value = 1365
$ws.Cells.Item(12, 6).Value = 61
$ws.Cells.Item(13, 6).Value = 2454
$ws.Cells.Item(14, 6).Value = 1161
$ws.Cells.Item(17, 6).Value = 196
$ws.Cells.Item(18, 6).Value = 27
$ws.Cells.Item(19, 6).Value = 58
$ws.Cells.Item(25, 6).Value = 262
$ws.Cells.Item(26, 6).Value = 37129
$ws.Cells.Item(27, 6).Value = 8
$ws.Cells.Item(33, 6).Value = 44
$ws.Cells.Item(35, 6).Value = 7
$ws.Cells.Item(38, 6).Value = 160
$ws.Cells.Item(42, 6).Value = 22
$ws.Cells.Item(43, 6).Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 859
$ws.Cells.Item(5, 6).Value = 237
$ws.Cells.Item(6, 6).Value = 2665
$ws.Cells.Item(7, 6).Value = 4434
$ws.Cells.Item(8, 6).Value = 110
$ws.Cells.Item(11, 6).Value = 561
$ws.Cells.Item(12, 6).Value = 379
$ws.Cells.Item(13, 6).Value = 85
$ws.Cells.Item(14, 6).Value = 515
$ws.Cells.Item(15, 6).Value = 164
$ws.Cells.Item(16, 6).Value = 170

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1824
$ws.Cells.Item(3, 6).Value = 859
$ws.Cells.Item(4, 6).Value = 237
$ws.Cells.Item(5, 6).Value = 4434
$ws.Cells.Item(6, 6).Value = 110
$ws.Cells.Item(7, 6).Value = 561
$ws.Cells.Item(9, 6).Value = 85
$ws.Cells.Item(10, 6).Value = 85
$ws.Cells.Item(11, 6).Value = 515
$ws.Cells.Item(12, 6).Value = 164
$ws.Cells.Item(13, 6).Value = 135
$ws.Cells.Item(14, 6).Value = 1072
$ws.Cells.Item(18, 6).Value = 1365
$ws.Cells.Item(19, 6).Value = 561
$ws.Cells.Item(20, 6).Value = 443
$ws.Cells.Item(21, 6).Value = 196
$ws.Cells.Item(22, 6).Value = 2454
$ws.Cells.Item(23, 6).Value = 1161
$ws.Cells.Item(24, 6).Value = 1347
$ws.Cells.Item(25, 6).Value = 1187
$ws.Cells.Item(26, 6).Value = 1387
$ws.Cells.Item(27, 6).Value = 86
$ws.Cells.Item(28, 6).Value = 196
$ws.Cells.Item(29, 6).Value = 58
$ws.Cells.Item(30, 6).Value = 1533
$ws.Cells.Item(31, 6).Value = 762
$ws.Cells.Item(32, 6).Value = 303
$ws.Cells.Item(34, 6).Value = 1141
$ws.Cells.Item(37, 6).Value = 986
$ws.Cells.Item(38, 6).Value = 262
$ws.Cells.Item(39, 6).Value = 927
$ws.Cells.Item(40, 6).Value = 873
$ws.Cells.Item(41, 6).Value = 814
$ws.Cells.Item(43, 6).Value = 1532
$ws.Cells.Item(44, 6).Value = 82
$ws.Cells.Item(45, 6).Value = 160
$ws.Cells.Item(47, 6).Value = 771
$ws.Cells.Item(48, 6).Value = 22
$ws.Cells.Item(49, 6).Value = 756
$ws.Cells.Item(50, 6).Value = 99
